$ws = $excel.ActiveWorkbook.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.489.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.835.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5329"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4053"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07569"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.86"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.54%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.113"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.320"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.57%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.642"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.43%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.000"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.14%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.88"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.11%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.838.21"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.24%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.74"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.03%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06595"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.67%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.59"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.78%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.089"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.73%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.514.43"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.34"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.96%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.115"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.68%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.458"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.52%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.09"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.73%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.59"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.80%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.048.54"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.46%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.85"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.93%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.132"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1095"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.88%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.681"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.03%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.656"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07162"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2266"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.32%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.265"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.69%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02342"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.73%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.867"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.29%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6296"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.48%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.35"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.196"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.19%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.413"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.61%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.49"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.714"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.81%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5866"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.85%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.02"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.994"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.30%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.194"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06915"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.47%  "
